# Update benchmark values on the "Image Streaming" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Image Streaming")

# Update the "Average Time Taken for sending each frame of image" column (D)
# with the new benchmark figures (values scaled up by 10x except the last one).
$ws.Range("D3").Value = "74 ms"
$ws.Range("D4").Value = "150 ms"
$ws.Range("D5").Value = "750 ms"
$ws.Range("D6").Value = "100 ms"
$ws.Range("D7").Value = "300 ms"
$ws.Range("D8").Value = "1.3 s"

# Update the selected cell in the sheet view.
$ws.Range("D12").Select()
